{"js": "// The tracked-changes \"_GoBack\" bookmark (Word's \"last edit position\" marker)\n// is moved from the end of the document (after the final tracked edits) back\n// to the empty paragraph right after the title - i.e. re-inserted at the top.\n//\n// Remove the stale bookmark first, then (re-)insert it around the second\n// paragraph of the document body (the blank line directly under the title).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[1];\nconst targetRange = targetParagraph.getRange(\"Whole\");\ntargetRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The tracked-changes \"_GoBack\" bookmark (Word's \"last edit position\" marker)\n# is moved from the end of the document (after the final tracked edits) back\n# to the empty paragraph right after the title - i.e. re-inserted at the top.\n$d = $word.ActiveDocument\n\n# Re-adding a bookmark with the same name (\"_GoBack\") automatically replaces\n# the existing one, so this both removes the stale bookmark at the end of the\n# document and creates the new one in a single step.\n$p = $d.Paragraphs(2)\n$d.Bookmarks.Add(\"_GoBack\", $p.Range)\n"}
